$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Ark1")

# Row 1 header values (column widths / n-counts)
$ws.Range("B1").Value = 16
$ws.Range("C1").Value = 20
$ws.Range("D1").Value = 16
$ws.Range("E1").Value = 20

# Row 2 (CON) updated values
$ws.Range("B2").Value = 60.750259275712338
$ws.Range("C2").Value = 60.124489566207195
$ws.Range("D2").Value = 53.771243284480605
$ws.Range("E2").Value = 67.871134820777954

# Row 3 (STR) updated values
$ws.Range("B3").Value = 31.805776781676283
$ws.Range("C3").Value = 56.603026870414674
$ws.Range("D3").Value = 56.955511943931079
$ws.Range("E3").Value = 83.044127267266163

# Update the active selection to match the edited range
$ws.Range("B1:E3").Select()
